# ngEO-WEBC internal issues sheet update: mark several issues as "Done"
# with assignees, merge the "popup widgets" description with its follow-up
# note, rename a couple of issue titles, and tidy the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Merge the continuation row (old row 5) into the "popup widgets" row (row 4) ---
$ws.Range("C4").Value = "not use jqm popups but a custom component to allow closing the widget.But if the JQM is upgraded the dismissable property is added so it can be used-to be decided"
$ws.Rows.Item(5).Delete()

# --- Row 4 now absorbs the extra row's height ---
$ws.Rows.Item(4).RowHeight = 90

# --- Rename a couple of issue titles/descriptions ---
$ws.Range("B5").Value = "Enter openSearch url for search"

# --- Status / Assignee updates ---
$ws.Range("D3").Value = "Done"
$ws.Range("E3").Value = "FL"

$ws.Range("D4").Value = "Done"
$ws.Range("E4").Value = "FL"

$ws.Range("D6").Value = "Done"

$ws.Range("D7").Value = "Done"
$ws.Range("E7").Value = "FL"

$ws.Range("E8").Value = "EM"

# --- Highlight the rows that are now Done with a green fill ---
$ws.Range("A3:E3").Interior.Color = 6604378
$ws.Range("A4:E4").Interior.Color = 6604378
$ws.Range("A6:E6").Interior.Color = 6604378
$ws.Range("A7:E7").Interior.Color = 6604378

# --- Move the active selection ---
$ws.Range("D2").Select() | Out-Null
